# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 11684
    "F3"  = 11304
    "F6"  = 1022
    "F9"  = 45
    "F11" = 10739
    "F12" = 4152
    "F19" = 126
    "F20" = 446
    "F21" = 11137
    "F22" = 10914
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
